# Apply updated vote-count figures after introducing a city-size
# (large/medium/small) criterion. Only specific cells on Sheet1 change;
# everything else in the workbook stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 132
$ws.Range("G3").Value = 122

$ws.Range("E5").Value = 489
$ws.Range("F5").Value = 279
$ws.Range("G5").Value = 270

$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 1023

$ws.Range("G8").Value = 12

$ws.Range("E9").Value = 40
$ws.Range("G9").Value = 81

$ws.Range("E10").Value = 355
$ws.Range("F10").Value = 311
$ws.Range("G10").Value = 334

$ws.Range("E12").Value = 385
$ws.Range("F12").Value = 273

$ws.Range("E14").Value = 42
$ws.Range("G14").Value = 27

$ws.Range("E15").Value = 16

$ws.Range("E16").Value = 559
$ws.Range("F16").Value = 480
$ws.Range("G16").Value = 506

$ws.Range("E18").Value = 56
$ws.Range("F18").Value = 81

$ws.Range("E21").Value = 16

$ws.Range("E22").Value = 11

$ws.Range("F23").Value = 436

$ws.Range("E24").Value = 58
$ws.Range("G24").Value = 87

$ws.Range("F25").Value = 490

$ws.Range("E26").Value = 776
$ws.Range("F26").Value = 688
$ws.Range("G26").Value = 800

$ws.Range("F27").Value = 23

$ws.Range("E30").Value = 562
$ws.Range("F30").Value = 644

$ws.Range("E31").Value = 420
$ws.Range("F31").Value = 303
$ws.Range("G31").Value = 266

$ws.Range("E33").Value = 75
$ws.Range("F33").Value = 102
